$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "linkage" (F) and "zygosity" (G) columns. Deleting the range
# shifts everything to the right of it left by two columns, matching the
# target layout (A:M instead of A:O).
$ws.Range("F1:G2").EntireColumn.Delete()
